$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 1: the empty paragraph right after "...already operating in
# four countries. https://selise.ch " (and right before "Software
# Engineer, Exabyting ...") gains run-level character formatting on
# its single empty run (bold / white / 11pt / dark-blue highlight).
# -----------------------------------------------------------------
$rngA = $d.Content.Duplicate
$foundA = $rngA.Find.Execute("already operating in four countries.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraA = $rngA.Paragraphs(1)
$emptyPara = $paraA.Next()
$emptyRng = $emptyPara.Range.Duplicate
$emptyRng.Collapse(1)
$emptyRng.Font.Bold = $true
$emptyRng.Font.BoldBi = $true
$emptyRng.Font.Color = 16777215
$emptyRng.Font.Size = 11
$emptyRng.Font.SizeBi = 11
$emptyRng.HighlightColorIndex = 9
$emptyRng.Font.Underline = 0

# -----------------------------------------------------------------
# Edit 2: split the run "TechServe4U " (trailing space) — that sits
# right after "Software Quality Assurance, " — into "TechServe4U"
# wrapped with a bookmark, plus a following run containing just the
# space.
# -----------------------------------------------------------------
$rngB = $d.Content.Duplicate
$foundB = $rngB.Find.Execute("TechServe4U ", $true, $false, $false, $false, $false, $true, 1, $false, "TechServe4U", 2)
$startB = $rngB.Start
$endB = $rngB.End
$afterB = $d.Range($endB, $endB)
$afterB.InsertAfter(" ")
$bmRangeB = $d.Range($startB, $endB)
$d.Bookmarks.Add("__DdeLink__420_1401661709", $bmRangeB)

# -----------------------------------------------------------------
# Edit 3: wrap the whole "TechServe4U, is situated at Michigan US..."
# paragraph with a bookmark (start right at the top of the paragraph,
# end right at the end of its text, before the paragraph mark), and
# fix the "perople" -> "people" typo inside it.
# -----------------------------------------------------------------
$rngC = $d.Content.Duplicate
$foundC = $rngC.Find.Execute("TechServe4U,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraC = $rngC.Paragraphs(1)
$startC = $paraC.Range.Start
$endC = $paraC.Range.End
$bmRangeC = $d.Range($startC, $endC - 1)
$d.Bookmarks.Add("__DdeLink__422_1401661709", $bmRangeC)

$rngD = $d.Content.Duplicate
$foundD = $rngD.Find.Execute("perople", $true, $false, $false, $false, $false, $true, 1, $false, "people", 2)

Write-Output ("edit1=" + $foundA + " edit2=" + $foundB + " edit3_bookmark=" + $foundC + " edit3_typo=" + $foundD)
